$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently spans A1:AC46. We add three new trailing columns
# (AD=Wins, AE=Losses, AF=Ties) holding the team's won/loss/tie record,
# repeated on every data row (2-46) since, per the commit message, this
# sheet keeps W/L/T alongside the per-player data instead of on a
# separate sheet.

# Copy the header formatting (bold font + border + centered alignment)
# from an existing header cell onto the three new header cells so they
# match the look of the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record for every data row.
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 87  # AD: Wins
    $ws.Cells.Item($r, 31).Value = 75  # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF: Ties
}
